# refactor: wrap table recognizer
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "4--10"
$ws.Range("E2").Value  = "10~9/L"

$ws.Range("C3").Value  = "'1"
$ws.Range("D3").Value  = "'2"

$ws.Range("C4").Value  = "'8"
$ws.Range("D4").Value  = "'0"

$ws.Range("C5").Value  = "'7.60"
$ws.Range("D5").Value  = "3--8"

$ws.Range("C6").Value  = "'0.20"
$ws.Range("D6").Value  = [char]0x2190 + "0.5--5"

$ws.Range("C7").Value  = "'0.10"
$ws.Range("D7").Value  = "0--1"

$ws.Range("A8").Value  = "ALY%"
$ws.Range("A9").Value  = "LIC%"

$ws.Range("D10").Value = [char]0x2193 + "2--7.7"
$ws.Range("E10").Value = "10.9/L"

$ws.Range("D11").Value = [char]0x2191 + "0.8--4"
$ws.Range("E11").Value = "10~9/L"

$ws.Range("A12").Value = "单核细胞"
$ws.Range("D12").Value = "0.12--0.8"

$ws.Range("A13").Value = "嗜酸性粒细胞"
$ws.Range("D13").Value = [char]0x2193 + "0.05--0.5"

$ws.Range("D14").Value = "0--0.1"

$ws.Range("C17").Value = "'126.00"
$ws.Range("F17").Value = "110--160"

$ws.Range("F18").Value = "3.5--5.5"

$ws.Range("F19").Value = "0.36--0.5"

$ws.Range("C20").Value = "'83.90"
$ws.Range("F20").Value = "0186--100"

$ws.Range("F21").Value = "26--31"

$ws.Range("C22").Value = "均"
$ws.Range("F22").Value = "310--370"

$ws.Range("C23").Value = "'41.40"
$ws.Range("F23").Value = "37--50"

$ws.Range("F24").Value = "11.5--14."

$ws.Range("C25").Value = "'5"
$ws.Range("F25").Value = "'4"

$ws.Range("F26").Value = "9--13"

$ws.Range("C27").Value = "'0.350"

$ws.Range("C28").Value = "'15.10"
$ws.Range("F28").Value = "9--17"
